# Apply the changes described by the diff:
# - PVThermal sheet: update parameters for "Type 1 -2.8m2" (A = 2.4, FR = 0.6, U_loss = 15)
#   and fix the "lamda_alpha" label typo to "lambda_alpha"; set active selection to B6.
# - BattModules sheet: set active selection to E12 (no data changes).
# - Hybrid OnGrid sheet: Idc_max_in changes from 25 to 45; set active selection to F14.

$wb = $excel.ActiveWorkbook

# --- PVThermal sheet ---
$wsPVThermal = $wb.Worksheets.Item("PVThermal")
$wsPVThermal.Activate()

# Fix typo in shared string label used by cell A7 ("lamda_alpha" -> "lambda_alpha")
$wsPVThermal.Range("A7").Value = "lambda_alpha"

# Update numeric parameters
$wsPVThermal.Range("B4").Value = 2.4
$wsPVThermal.Range("B5").Value = 0.6
$wsPVThermal.Range("B6").Value = 15

$wsPVThermal.Range("B6").Select()

# --- BattModules sheet ---
$wsBattModules = $wb.Worksheets.Item("BattModules")
$wsBattModules.Activate()
$wsBattModules.Range("E12").Select()

# --- Hybrid OnGrid sheet ---
$wsHybrid = $wb.Worksheets.Item("Hybrid OnGrid")
$wsHybrid.Activate()
$wsHybrid.Range("B10").Value = 45
$wsHybrid.Range("F14").Select()

# Restore PVThermal as the active sheet/tab
$wsPVThermal.Activate()
$wsPVThermal.Range("B6").Select()
